$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# 1. Insert a new row at 53 (everything from the old row 53 downward
#    shifts down by one row: old 53->54, 54->55, ... 60->61).
# ------------------------------------------------------------------
$ws.Rows.Item(53).Insert()

# Re-use the formatting of an existing "monthly data" row (row 22 has
# the exact same style pattern we need: s=11,12,12,13,13,13,14,12,12,13,13,14,13)
# so that the newly inserted row gets the correct (already-existing)
# style indexes instead of brand new ones.
$ws.Range("A22:M22").Copy()
$ws.Range("A53:M53").PasteSpecial(-4122)   # xlPasteFormats
$excel.CutCopyMode = $false

# ------------------------------------------------------------------
# 2. Fill in the new "November" monthly data row (row 53).
# ------------------------------------------------------------------
$ws.Range("A53").Value2 = "November"
$ws.Range("B53").Value2 = 0
$ws.Range("C53").Value2 = 0
$ws.Range("D53").Value2 = "--"
$ws.Range("E53").Value2 = "--"
$ws.Range("F53").Value2 = "--"
$ws.Range("G53").Value2 = 0
$ws.Range("H53").Value2 = 613
$ws.Range("I53").Value2 = 589
$ws.Range("J53").Value2 = "W"
$ws.Range("K53").Value2 = "W"
$ws.Range("L53").Value2 = 6.5
$ws.Range("M53").Value2 = "W"

# ------------------------------------------------------------------
# 3. Update the "Year to Date" section (now rows 55-57, was 54-56)
#    with the refreshed figures.
# ------------------------------------------------------------------
# Year 2014 (row 55)
$ws.Range("H55").Value2 = 5177
$ws.Range("I55").Value2 = 5130
$ws.Range("L55").Value2 = 4.7

# Year 2015 (row 56)
$ws.Range("H56").Value2 = 5975
$ws.Range("I56").Value2 = 5863
$ws.Range("L56").Value2 = 5.5

# Year 2016 (row 57)
$ws.Range("H57").Value2 = 7438
$ws.Range("I57").Value2 = 7217
$ws.Range("L57").Value2 = 6.7

# ------------------------------------------------------------------
# 4. Update the "Rolling 12 Months Ending in ..." section
#    (now rows 59-60, was 58-59).
# ------------------------------------------------------------------
# Year 2015 (row 59)
$ws.Range("H59").Value2 = 6647
$ws.Range("I59").Value2 = 6529
$ws.Range("L59").Value2 = 5.6

# Year 2016 (row 60)
$ws.Range("H60").Value2 = 7962
$ws.Range("I60").Value2 = 7724
$ws.Range("L60").Value2 = 6.6

# ------------------------------------------------------------------
# 5. Text updates: title, and the "Rolling 12 Months Ending in" label.
# ------------------------------------------------------------------
$ws.Range("A1").Value2 = "Table 4.4. Receipts, Average Cost, and Quality of Fossil Fuels: Commerical Sector, 2006 - November 2016 (continued)"
$ws.Range("A58").Value2 = "Rolling 12 Months Ending in November"
